$wb = $excel.ActiveWorkbook

# --- weibull ---
$ws = $wb.Worksheets.Item("weibull")
$ws.Range("B2").Value = -2.8473150791631
$ws.Range("C2").Value = 0.144695958439284
$ws.Range("B3").Value = -0.0750160783595377
$ws.Range("C3").Value = 0.0858304841970044

# --- lognormal ---
$ws = $wb.Worksheets.Item("lognormal")
$ws.Range("B2").Value = 2.25444620048039
$ws.Range("C2").Value = 0.222341048830984
$ws.Range("B3").Value = -0.911033282332893
$ws.Range("C3").Value = 0.102491288079624

# --- llogis ---
$ws = $wb.Worksheets.Item("llogis")
$ws.Range("B2").Value = -2.4143783123601
$ws.Range("C2").Value = 0.137942025672175
$ws.Range("B3").Value = 0.462959154088631
$ws.Range("C3").Value = 0.128537936727382

# --- gompertz ---
$ws = $wb.Worksheets.Item("gompertz")
$ws.Range("B2").Value = -2.51123869942539
$ws.Range("C2").Value = 0.121563063431743
$ws.Range("B3").Value = -0.0369299932610897
$ws.Range("C3").Value = 0.010499509130997

# --- weibull cov ---
$ws = $wb.Worksheets.Item("weibull cov")
$ws.Range("A2").Value = 0.0209369203886629
$ws.Range("B2").Value = -0.00586329307363494
$ws.Range("A3").Value = -0.00586329307363494
$ws.Range("B3").Value = 0.00736687201749222

# --- lognormal cov ---
$ws = $wb.Worksheets.Item("lognormal cov")
$ws.Range("A2").Value = 0.0494355419952621
$ws.Range("B2").Value = -0.0192375999914305
$ws.Range("A3").Value = -0.0192375999914305
$ws.Range("B3").Value = 0.0105044641322206

# --- llogis cov ---
$ws = $wb.Worksheets.Item("llogis cov")
$ws.Range("A2").Value = 0.0190280024465431
$ws.Range("B2").Value = 0.00808717390768899
$ws.Range("A3").Value = 0.00808717390768899
$ws.Range("B3").Value = 0.0165220011781323

# --- gompertz cov ---
$ws = $wb.Worksheets.Item("gompertz cov")
$ws.Range("A2").Value = 0.01477757839091
$ws.Range("B2").Value = -0.000462489919387872
$ws.Range("A3").Value = -0.000462489919387872
$ws.Range("B3").Value = 0.000110239691991889
